$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# Insert a single new leave-card row right after row 94 (i.e. before
# the old row 95). Because row 95 sits in the middle of Table1, this
# single insert cascades all the way to the bottom of the table:
#   old row 95..133  -> new row 96..134   (ordinary data rows)
#   old row 134       -> new row 135       (the special last/totals-
#                                            styled row)
# and the freshly-opened row 95 becomes the *new* entry we still need
# to fill in (it does NOT create a second phantom row at the bottom).
# ------------------------------------------------------------------
$ws.Rows.Item(95).Insert()

# Copy the formatting of row 93 (an existing "SL(1-0-00)" entry row)
# onto the freshly inserted row 95 so every cell gets the same styles
# used by sibling rows (A=40,B=20,C=13,D=39,E=9,F=20,G=13,H=39,I=9,J=11,K=48).
$ws.Range("A93:K93").Copy() | Out-Null
$ws.Range("A95:K95").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 94: keep its existing PERIOD date, add the SL(1-0-00) leave entry.
$ws.Cells.Item(94, 2).Value2 = "SL(1-0-00)"
$ws.Cells.Item(94, 8).Value2 = 1
$ws.Cells.Item(94, 11).Value2 = 45016

# Row 95 (new row): no PERIOD date, same SL(1-0-00) leave entry.
$ws.Cells.Item(95, 1).ClearContents() | Out-Null
$ws.Cells.Item(95, 2).Value2 = "SL(1-0-00)"
$ws.Cells.Item(95, 8).Value2 = 1
$ws.Cells.Item(95, 11).Value2 = 45027

# Restore the calculated-column formula on both rows (G column) plus
# the row at the very bottom of the table (row 135) whose structured
# reference sometimes gets rewritten by the row-shift.
$calcFormula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Cells.Item(94, 7).Formula = $calcFormula
$ws.Cells.Item(95, 7).Formula = $calcFormula
$ws.Cells.Item(135, 7).Formula = $calcFormula

# Grow the table definition to match the new bottom row (135).
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K135")) | Out-Null

# ------------------------------------------------------------------
# Refresh the cached selection so it mirrors the saved file.
# ------------------------------------------------------------------
$ws.Range("B96").Select() | Out-Null

$wb.Application.Calculate()
